$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in H1, using the same style as the existing
# header cells (bold/border/centered, style index 1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Add the new Save value in H2 (plain number, no special style, like F2/G2)
$ws.Range("H2").Value = 0
